$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1260
$ws.Range("I4").Value = 325
$ws.Range("K4").Value = 325
$ws.Range("M4").Value = -211
$ws.Range("H6").Value = 249
$ws.Range("I6").Value = 249
$ws.Range("K6").Value = 747
$ws.Range("M6").Value = -635
$ws.Range("H18").Value = 138.16667
$ws.Range("I18").Value = 138.16667
$ws.Range("K18").Value = 138.16667
$ws.Range("M18").Value = 145.83333
$ws.Range("H96").Value = 721.4
$ws.Range("I96").Value = 504
$ws.Range("J96").Value = 866.3333
$ws.Range("K96").Value = 1512
$ws.Range("L96").Value = 2598.9999
$ws.Range("M96").Value = -139
$ws.Range("N96").Value = -5344.9999
$ws.Range("H132").Value = 4430.359
$ws.Range("I132").Value = 4451.625
$ws.Range("K132").Value = 13354.875
$ws.Range("M132").Value = -10824.875
$ws.Range("H134").Value = 42996.688
$ws.Range("J134").Value = 42996.688
$ws.Range("L134").Value = 42996.688
$ws.Range("N134").Value = -53136.688
$ws.Range("H137").Value = 1334.0476
$ws.Range("I137").Value = 862.6
$ws.Range("K137").Value = 2587.8
$ws.Range("M137").Value = -37.80000000000018

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2354.9246
$ws.Range("I32").Value = 2192.5293
$ws.Range("K32").Value = 2192.5293
$ws.Range("M32").Value = -1905.5293
$ws.Range("H45").Value = 1660.75
$ws.Range("I45").Value = 1159.5
$ws.Range("K45").Value = 1159.5
$ws.Range("M45").Value = -782.5
$ws.Range("H97").Value = 11550.833
$ws.Range("I97").Value = 15098.25
$ws.Range("K97").Value = 15098.25
$ws.Range("M97").Value = -14602.25
$ws.Range("H122").Value = 3254.4443
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 66835.60000000001
$ws.Range("J132").Value = 66835.60000000001
$ws.Range("L132").Value = 66835.60000000001
$ws.Range("N132").Value = -76955.60000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3767.15
$ws.Range("I31").Value = 1999
$ws.Range("J31").Value = 3860.2104
$ws.Range("K31").Value = 1999
$ws.Range("L31").Value = 3860.2104
$ws.Range("M31").Value = -1704
$ws.Range("N31").Value = -4450.2104
$ws.Range("H34").Value = 3767.15
$ws.Range("I34").Value = 1999
$ws.Range("J34").Value = 3860.2104
$ws.Range("K34").Value = 1999
$ws.Range("L34").Value = 3860.2104
$ws.Range("M34").Value = -1797
$ws.Range("N34").Value = -4264.2104
$ws.Range("H58").Value = 4501.6665
$ws.Range("I58").Value = 3543.7856
$ws.Range("J58").Value = 6417.4287
$ws.Range("K58").Value = 3543.7856
$ws.Range("L58").Value = 6417.4287
$ws.Range("M58").Value = -3340.7856
$ws.Range("N58").Value = -6823.4287
$ws.Range("H94").Value = 2143.7334
$ws.Range("J94").Value = 2305.75
$ws.Range("L94").Value = 2305.75
$ws.Range("N94").Value = -3207.75
$ws.Range("H136").Value = 4501.6665
$ws.Range("I136").Value = 3543.7856
$ws.Range("J136").Value = 6417.4287
$ws.Range("K136").Value = 10631.3568
$ws.Range("L136").Value = 19252.2861
$ws.Range("M136").Value = -8081.356800000001
$ws.Range("N136").Value = -24352.2861

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 569.2727
$ws.Range("J12").Value = 626.1
$ws.Range("L12").Value = 1878.3
$ws.Range("N12").Value = -2224.3
$ws.Range("H81").Value = 8999.799999999999
$ws.Range("I81").Value = 8999.799999999999
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 26999.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -25876.4
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 8999.799999999999
$ws.Range("I84").Value = 8999.799999999999
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 80998.2
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -75382.2
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 595
$ws.Range("I113").Value = 521.93335
$ws.Range("J113").Value = 716.7778
$ws.Range("K113").Value = 1565.80005
$ws.Range("L113").Value = 2150.3334
$ws.Range("M113").Value = 604.1999499999999
$ws.Range("N113").Value = -6490.3334
$ws.Range("H117").Value = 2123.7058
$ws.Range("J117").Value = 2695
$ws.Range("L117").Value = 8085
$ws.Range("N117").Value = -14969
$ws.Range("H121").Value = 4200
$ws.Range("I121").Value = 2950
$ws.Range("J121").Value = 4366.6665
$ws.Range("K121").Value = 8850
$ws.Range("L121").Value = 13099.9995
$ws.Range("M121").Value = -7540
$ws.Range("N121").Value = -15719.9995
$ws.Range("H134").Value = 9656.275
$ws.Range("I134").Value = 4119.5293
$ws.Range("K134").Value = 12358.5879
$ws.Range("M134").Value = -7288.5879

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2783.25
$ws.Range("J80").Value = 2600
$ws.Range("L80").Value = 2600
$ws.Range("N80").Value = -4596
$ws.Range("H83").Value = 2783.25
$ws.Range("J83").Value = 2600
$ws.Range("L83").Value = 13000
$ws.Range("N83").Value = -22984
$ws.Range("H92").Value = 10566.571
$ws.Range("I92").Value = 928
$ws.Range("J92").Value = 12173
$ws.Range("K92").Value = 928
$ws.Range("L92").Value = 12173
$ws.Range("M92").Value = 944
$ws.Range("N92").Value = -15917
$ws.Range("H122").Value = 3482.8333
$ws.Range("J122").Value = 4666
$ws.Range("L122").Value = 13998
$ws.Range("N122").Value = -18898

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1327.8572
$ws.Range("I22").Value = 565.6667
$ws.Range("J22").Value = 1899.5
$ws.Range("K22").Value = 565.6667
$ws.Range("L22").Value = 1899.5
$ws.Range("M22").Value = -270.6667
$ws.Range("N22").Value = -2489.5
$ws.Range("H27").Value = 1327.8572
$ws.Range("I27").Value = 565.6667
$ws.Range("J27").Value = 1899.5
$ws.Range("K27").Value = 565.6667
$ws.Range("L27").Value = 1899.5
$ws.Range("M27").Value = -458.6667
$ws.Range("N27").Value = -2113.5
$ws.Range("H46").Value = 2278.8333
$ws.Range("J46").Value = 2534.6
$ws.Range("L46").Value = 2534.6
$ws.Range("N46").Value = -2910.6
$ws.Range("H68").Value = 2532.9375
$ws.Range("I68").Value = 2457
$ws.Range("K68").Value = 2457
$ws.Range("M68").Value = -1708
$ws.Range("H71").Value = 2532.9375
$ws.Range("I71").Value = 2457
$ws.Range("K71").Value = 12285
$ws.Range("M71").Value = -8541
$ws.Range("H82").Value = 5305.8667
$ws.Range("I82").Value = 2450
$ws.Range("J82").Value = 7209.778
$ws.Range("K82").Value = 2450
$ws.Range("L82").Value = 7209.778
$ws.Range("M82").Value = -2089
$ws.Range("N82").Value = -7931.778
$ws.Range("H85").Value = 5305.8667
$ws.Range("I85").Value = 2450
$ws.Range("J85").Value = 7209.778
$ws.Range("K85").Value = 2450
$ws.Range("L85").Value = 7209.778
$ws.Range("M85").Value = -1202
$ws.Range("N85").Value = -9705.778
$ws.Range("H132").Value = 5417.9375
$ws.Range("I132").Value = 4733.1665
$ws.Range("K132").Value = 14199.4995
$ws.Range("M132").Value = -11669.4995

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9242.571
$ws.Range("I62").Value = 8766.666999999999
$ws.Range("J62").Value = 9599.5
$ws.Range("K62").Value = 8766.666999999999
$ws.Range("L62").Value = 9599.5
$ws.Range("M62").Value = -8142.666999999999
$ws.Range("N62").Value = -10847.5
$ws.Range("H65").Value = 9242.571
$ws.Range("I65").Value = 8766.666999999999
$ws.Range("J65").Value = 9599.5
$ws.Range("K65").Value = 43833.335
$ws.Range("L65").Value = 47997.5
$ws.Range("M65").Value = -40713.335
$ws.Range("N65").Value = -54237.5
$ws.Range("H113").Value = 450
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 450
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -5690
